# Apply updates described by the diff (output generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: 展览
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value = 514
$ws.Range("F4").Value = 491
$ws.Range("F7").Value = 942
$ws.Range("F8").Value = 739
$ws.Range("F9").Value = 187
$ws.Range("F10").Value = 48

$ws.Range("C11").Value = "广州·ANM国际学院&国际集团管培生培养计划沉浸式国乙ONLY（取消）"
$ws.Range("G11").Value = "不可售"

$ws.Range("F12").Value = 770
$ws.Range("F16").Value = 1292
$ws.Range("F17").Value = 113
$ws.Range("F19").Value = 1088
$ws.Range("F21").Value = 1285
$ws.Range("F22").Value = 650
$ws.Range("F24").Value = 1243
$ws.Range("F28").Value = 986
$ws.Range("F29").Value = 18
$ws.Range("F31").Value = 1320

# ---------------------------------------------------------------
# Sheet: 演出
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

$ws.Range("G3").Value = 888
$ws.Range("F7").Value = 7

# ---------------------------------------------------------------
# Sheet: 本地生活
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")

$ws.Range("F2").Value = 719

# ---------------------------------------------------------------
# Sheet: 全部类型
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F2").Value = 719
$ws.Range("F3").Value = 514
$ws.Range("F6").Value = 491
$ws.Range("G7").Value = 888
$ws.Range("G8").Value = 888
$ws.Range("F14").Value = 7
$ws.Range("F15").Value = 942
$ws.Range("F16").Value = 739
$ws.Range("F17").Value = 187
$ws.Range("F19").Value = 48

$ws.Range("C23").Value = "广州·ANM国际学院&国际集团管培生培养计划沉浸式国乙ONLY（取消）"
$ws.Range("G23").Value = "不可售"

$ws.Range("F25").Value = 770
$ws.Range("F29").Value = 1292
$ws.Range("F30").Value = 113
$ws.Range("F32").Value = 1088
$ws.Range("F34").Value = 1285
$ws.Range("F35").Value = 650
$ws.Range("F37").Value = 1243
$ws.Range("F43").Value = 986
$ws.Range("F44").Value = 18
$ws.Range("F46").Value = 1320
